$wb = $excel.ActiveWorkbook

# Add the new worksheet as the last tab and name it to match the
# 20080516_Han-et-al paper.
$lastIndex = $wb.Worksheets.Count
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$newSheet.Name = "20080516_Han-et-al"

# rsid header + the 38 SNP ids from the paper, in sheet order.
$values = @(
    "rsid",
    "rs12913832",
    "rs1667394",
    "rs12203592",
    "rs258322",
    "rs4785763",
    "rs6497268",
    "rs8039195",
    "rs11855019",
    "rs11636232",
    "rs8049897",
    "rs4238833",
    "rs4408545",
    "rs7204478",
    "rs4904866",
    "rs12896399",
    "rs7174027",
    "rs7183877",
    "rs7196459",
    "rs164741",
    "rs7188458",
    "rs8033165",
    "rs35391",
    "rs7495174",
    "rs1635168",
    "rs8007923",
    "rs10861741",
    "rs28777",
    "rs9806558",
    "rs9392056",
    "rs4778211",
    "rs2493040",
    "rs6918152",
    "rs2353033",
    "rs12142165",
    "rs7195066",
    "rs2241039",
    "rs8028689",
    "rs16950987"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match the saved selection/active cell state on the new sheet.
[void]$newSheet.Range("A40").Select()
